$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "딱 내가 배우고 싶은 것만 배우면 되는데 뭐가 이렇게 많아?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/why-pre-requisites/#utm_source=rss&utm_medium=rss&utm_campaign=why-pre-requisites"

$ws.Range("D28").Value = "RRT 경로 생성 알고리즘"
$ws.Range("E28").Value = "https://ropiens.tistory.com/192"

$ws.Range("D32").Value = "분산분석 (ANOVA : ANalysis Of VAriance)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/379"

$ws.Range("D36").Value = "Introduction to Graph Neural Networks(Spectral Graph Convolution)"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/371"

$ws.Range("D51").Value = "[flask] ubuntu에 flask_mysqldb 설치 중 에러 해결 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/flask-ubuntu%EC%97%90-flaskmysqldb-%EC%84%A4%EC%B9%98-%EC%A4%91-%EC%97%90%EB%9F%AC-%ED%95%B4%EA%B2%B0-%EB%B0%A9%EB%B2%95"

$ws.Range("D52").Value = "숨은 DS"

$wb.Save()
